$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = "Time elapsed calculated by wrapping Logstash with the Linux time command"
$ws.Range("C22").Value = "Example:"
$ws.Range("C24").Value = "time sudo bin/logstash -w 6 -f /path/to/logstash_config_file.conf"

$ws.Range("C25").Select()
